# Apply edits described by the diff:
# 1. Update shared-string text for two header cells (append "*")
# 2. Move the active selection from K2 to G2
# 3. Add a list data-validation ("yes,no") on F2:G2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update header text -------------------------------------------------
$ws.Range("F1").Value = "Needs Delivery Note ?*"
$ws.Range("G1").Value = "Has Attchment ?*"

# --- 2. Move selection to G2 -------------------------------------------------
$ws.Range("G2").Select()

# --- 3. Data validation (list: yes,no) on F2:G2 -----------------------------
$rng = $ws.Range("F2:G2")
$rng.Validation.Delete()
$rng.Validation.Add(3, 1, 1, '"yes,no"')
